# Changed figure labels and units
# - Rename the "excretion rate N" / "excretion rate P" trait-bucket labels
#   (col B, rows 36-39) to the more descriptive "excretion rate of ammonia"
#   and "excretion rate of phosphate".
# - Update the active view/selection to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B36").Value = "excretion rate of ammonia"
$ws.Range("B37").Value = "excretion rate of ammonia"
$ws.Range("B38").Value = "excretion rate of phosphate"
$ws.Range("B39").Value = "excretion rate of phosphate"

# Reflect the final cursor position / scroll of the sheet, matching the
# saved view state (top-left visible cell A22, active cell B39).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("B39").Select()
